$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 45456372
$ws.Range("I12").Value = 2601.818
$ws.Range("J12").Value = 90910140
$ws.Range("K12").Value = 2601.818
$ws.Range("L12").Value = 90910140
$ws.Range("M12").Value = -2431.818
$ws.Range("N12").Value = -90910480
$ws.Range("H16").Value = 17524.555
$ws.Range("I16").Value = 2975.25
$ws.Range("J16").Value = 29164
$ws.Range("K16").Value = 2975.25
$ws.Range("L16").Value = 29164
$ws.Range("M16").Value = -2745.25
$ws.Range("N16").Value = -29624
$ws.Range("H21").Value = 16267.934
$ws.Range("I21").Value = 18431.285
$ws.Range("K21").Value = 18431.285
$ws.Range("M21").Value = -17963.285
$ws.Range("H23").Value = 16267.934
$ws.Range("I23").Value = 18431.285
$ws.Range("K23").Value = 18431.285
$ws.Range("M23").Value = -18197.285
$ws.Range("H29").Value = 559.7143
$ws.Range("I29").Value = 123.6
$ws.Range("K29").Value = 370.8
$ws.Range("M29").Value = -89.79999999999995
$ws.Range("H38").Value = 394.14285
$ws.Range("I38").Value = 394.14285
$ws.Range("K38").Value = 1182.42855
$ws.Range("M38").Value = -810.4285500000001
$ws.Range("H76").Value = 3178.762
$ws.Range("J76").Value = 3259
$ws.Range("L76").Value = 3259
$ws.Range("N76").Value = -3889
$ws.Range("H79").Value = 3178.762
$ws.Range("J79").Value = 3259
$ws.Range("L79").Value = 3259
$ws.Range("N79").Value = -5443
$ws.Range("H112").Value = 1332.7273
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1332.7273
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3998.1819
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -6214.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1675.7931
$ws.Range("I45").Value = 1762.8235
$ws.Range("J45").Value = 1552.5
$ws.Range("K45").Value = 1762.8235
$ws.Range("L45").Value = 1552.5
$ws.Range("M45").Value = -1385.8235
$ws.Range("N45").Value = -2306.5
$ws.Range("H122").Value = 2263.7334
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2263.7334
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6791.2002
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -11691.2002
$ws.Range("H139").Value = 35184.25
$ws.Range("J139").Value = 35184.25
$ws.Range("L139").Value = 35184.25
$ws.Range("N139").Value = -45464.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29726.666
$ws.Range("J55").Value = 29726.666
$ws.Range("L55").Value = 29726.666
$ws.Range("N55").Value = -30272.666
$ws.Range("H105").Value = 2242.4167
$ws.Range("I105").Value = 2001.5
$ws.Range("J105").Value = 2483.3333
$ws.Range("K105").Value = 2001.5
$ws.Range("L105").Value = 2483.3333
$ws.Range("M105").Value = -254.5
$ws.Range("N105").Value = -5977.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 28000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 28000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 28000
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -28630
$ws.Range("H122").Value = 1297.5385
$ws.Range("I122").Value = 1350
$ws.Range("J122").Value = 1288
$ws.Range("K122").Value = 4050
$ws.Range("L122").Value = 3864
$ws.Range("M122").Value = -1600
$ws.Range("N122").Value = -8764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 8143
$ws.Range("J105").Value = 8143
$ws.Range("L105").Value = 24429
$ws.Range("N105").Value = -29671
$ws.Range("H121").Value = 1468
$ws.Range("I121").Value = 625
$ws.Range("J121").Value = 1655.3334
$ws.Range("K121").Value = 1875
$ws.Range("L121").Value = 4966.0002
$ws.Range("M121").Value = -565
$ws.Range("N121").Value = -7586.0002
$ws.Range("H131").Value = 829.2062
$ws.Range("I131").Value = 456.58334
$ws.Range("J131").Value = 881.81177
$ws.Range("K131").Value = 1369.75002
$ws.Range("L131").Value = 2645.43531
$ws.Range("M131").Value = 3670.24998
$ws.Range("N131").Value = -12725.43531
$ws.Range("H132").Value = 2479.3428
$ws.Range("I132").Value = 2313
$ws.Range("J132").Value = 2842.2727
$ws.Range("K132").Value = 20817
$ws.Range("L132").Value = 25580.4543
$ws.Range("M132").Value = -18287
$ws.Range("N132").Value = -30640.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2534.1052
$ws.Range("I80").Value = 2534.1052
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2534.1052
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1536.1052
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 2534.1052
$ws.Range("I83").Value = 2534.1052
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12670.526
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7678.526
$ws.Range("N83").Value = ""
$ws.Range("H101").Value = 29966.666
$ws.Range("J101").Value = 29966.666
$ws.Range("L101").Value = 29966.666
$ws.Range("N101").Value = -36456.666
$ws.Range("H132").Value = 3221.1
$ws.Range("I132").Value = 3188.2666
$ws.Range("J132").Value = 3319.6
$ws.Range("K132").Value = 9564.799800000001
$ws.Range("L132").Value = 9958.799999999999
$ws.Range("M132").Value = -7034.799800000001
$ws.Range("N132").Value = -15018.8
$ws.Range("H133").Value = 34754
$ws.Range("J133").Value = 34754
$ws.Range("L133").Value = 34754
$ws.Range("N133").Value = -44874
$ws.Range("H137").Value = 56780
$ws.Range("J137").Value = 56780
$ws.Range("L137").Value = 56780
$ws.Range("N137").Value = -66980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4831.7896
$ws.Range("I122").Value = 5799.9165
$ws.Range("J122").Value = 3172.1428
$ws.Range("K122").Value = 17399.7495
$ws.Range("L122").Value = 9516.428400000001
$ws.Range("M122").Value = -14949.7495
$ws.Range("N122").Value = -14416.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6980562
$ws.Range("I122").Value = 10872357
$ws.Range("J122").Value = 95077.92
$ws.Range("K122").Value = 32617071
$ws.Range("L122").Value = 285233.76
$ws.Range("M122").Value = -32614621
$ws.Range("N122").Value = -290133.76
